$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 84.81481481481481
$ws.Range("A3").Value = 83.703703703703695
$ws.Range("A4").Value = 85.18518518518519
$ws.Range("A5").Value = 79.629629629629633
$ws.Range("A6").Value = 81.481481481481481
$ws.Range("A7").Value = 81.481481481481481
$ws.Range("A8").Value = 88.888888888888886
$ws.Range("A9").Value = 88.888888888888886
$ws.Range("A10").Value = 88.888888888888886
$ws.Range("A11").Value = 88.518518518518519
$ws.Range("A12").Value = 76.296296296296291
$ws.Range("A13").Value = 83.333333333333343
$ws.Range("A14").Value = 88.518518518518519
$ws.Range("A15").Value = 87.407407407407405
$ws.Range("A16").Value = 88.148148148148152
$ws.Range("A17").Value = 80
$ws.Range("A18").Value = 81.851851851851848
$ws.Range("A19").Value = 84.074074074074076
$ws.Range("A20").Value = 85.925925925925924
$ws.Range("A21").Value = 87.407407407407405
$ws.Range("A22").Value = 87.407407407407405
$ws.Range("A23").Value = 77.777777777777786
$ws.Range("A24").Value = 77.407407407407405
$ws.Range("A25").Value = 76.296296296296291
$ws.Range("A26").Value = 87.407407407407405
$ws.Range("A27").Value = 84.444444444444443
$ws.Range("A28").Value = 85.555555555555557
$ws.Range("A29").Value = 83.333333333333343
$ws.Range("A30").Value = 81.851851851851848
$ws.Range("A31").Value = 83.703703703703695
$ws.Range("A32").Value = 81.851851851851848
$ws.Range("A33").Value = 82.962962962962962
$ws.Range("A34").Value = 82.222222222222214
$ws.Range("A35").Value = 80
$ws.Range("A36").Value = 80.370370370370367
$ws.Range("A37").Value = 78.518518518518519
$ws.Range("A38").Value = 80.370370370370367
$ws.Range("A39").Value = 77.407407407407405
$ws.Range("A40").Value = 77.407407407407405
$ws.Range("A41").Value = 85.925925925925924
$ws.Range("A42").Value = 87.037037037037038
$ws.Range("A43").Value = 86.296296296296291
$ws.Range("A44").Value = 85.555555555555557
$ws.Range("A45").Value = 85.925925925925924
$ws.Range("A46").Value = 85.555555555555557
$ws.Range("A47").Value = 78.888888888888886
$ws.Range("A48").Value = 79.629629629629633
$ws.Range("A49").Value = 84.074074074074076
